$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Data edits on "vaccines" sheet ---

# Header: "supply chain" -> "cooling temperature"
$ws1.Range("K1").Value = "cooling temperature"

# Fix typos "gab" -> "gap" in vaccination plan column
$ws1.Range("I3").Value = "4 weeks gap between doses"
$ws1.Range("I2").Value = "3 weeks gap between shots"

# Clear the "(link)" placeholders in the type info column for AstraZeneca / Janssen-Cilag,
# keeping their highlighted (yellow) formatting
$ws1.Range("F4").Value = ""

$ws1.Range("F5").Interior.Color = 65535
$ws1.Range("F5").Value = ""

# Clear the reserved-doses placeholder for Janssen-Cilag
$ws1.Range("N5").Value = ""

# --- Active sheet / selection changes ---
# Make "vaccines" the active sheet with N13 selected
$ws1.Activate()
$ws1.Range("N13").Select()

# Restore the selection on "vaccines links" to its original cell
$ws2.Range("B9").Select()
